$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Working Papers")

# Update the VSL paper abstract in C3 first (so shared-string ordering matches source)
$ws.Range("C3").Value = 'Economists often study non-market goods such as health and air quality. This paper introduces a new method to estimate demand for such amenities and applies it to measure the value of a statistical life (VSL) in Kenya. My approach is to update beliefs about the life-saving efficacy of a product (a motorcycle helmet) and elicit product choice. This generates instruments allowing one to use subjective beliefs to estimate demand, rather than assuming rational expectations. This method does not require beliefs to be reported error-free but does require classical mismeasurement. I validate this assumption using features of the experimental design. The estimated VSL is \$224, near the left tail of Kenyan estimates. Standard methods for estimating VSL produce skewed results, driven by severe violations of rational expectations. These findings help explain low observed demand for many health products and suggest that directing more development aid towards consumption may increase welfare. '

# Update the Job Market Paper (JMP) abstract in C2
$ws.Range("C2").Value = 'Economists typically assume that firms are risk neutral. Yet many enterprises, especially in low and middle-income economies, are small and owner-operated, making household consumption sensitive to business risk. This may cause the risk preferences of firm owners to affect firm decisions. This paper demonstrates that small retailers in Kenya are risk averse, leading them to under-adopt a new product because they face uncertain demand. I develop a model in which risk averse firms learn about demand through stocking decisions, then test its predictions using two field experiments. The first establishes that risk aversion affects the stocking decisions of enterprises. I construct a new test for risk aversion by offering treated firms an insurance contract that lowers expected profits from a new product and reduces the risk of losses, leading to a 50\% increase in adoption which implies a rejection of risk neutrality. The second experiment shows that \textit{temporarily} reducing inventory risk leads firms to \textit{permanently} stock a profitable new product because they overcome demand uncertainty through learning. These results suggest a departure from perfect competition, since risk averse firms are not competed out of the market. I document information externalities consistent with this prediction: enterprises that are nearby randomly induced entrants stock the new product at a higher rate. These results suggest that risk aversion can impede product diffusion and firm growth, challenging the common assumption that small firms are risk neutral.'

# Add a new value in E3 (flag column), matching E2
$ws.Range("E3").Value = 1

# The VSL abstract is now shorter, so row 3 auto-fits to a smaller height
$ws.Rows.Item(3).RowHeight = 348

Write-Host "Edit complete"
